# Auto-generated edit script applying market-price refresh changes
# to the Leviathan_Profits workbook (per sheet: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 346822.38
$ws.Range("J17").Value = 346822.38
$ws.Range("L17").Value = 1040467.14
$ws.Range("N17").Value = -1040803.14
$ws.Range("H74").Value = 4422.769
$ws.Range("I74").Value = 3356.5715
$ws.Range("K74").Value = 3356.5715
$ws.Range("M74").Value = -2420.5715
$ws.Range("H77").Value = 4422.769
$ws.Range("I77").Value = 3356.5715
$ws.Range("K77").Value = 16782.8575
$ws.Range("M77").Value = -12102.8575
$ws.Range("H97").Value = 1542.4286
$ws.Range("J97").Value = 1542.4286
$ws.Range("L97").Value = 4627.2858
$ws.Range("N97").Value = -5619.2858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4308.0356
$ws.Range("I32").Value = 3833.5095
$ws.Range("J32").Value = 12691.333
$ws.Range("K32").Value = 3833.5095
$ws.Range("L32").Value = 12691.333
$ws.Range("M32").Value = -3546.5095
$ws.Range("N32").Value = -13265.333
$ws.Range("H74").Value = 1475.3334
$ws.Range("I74").Value = 1366.1177
$ws.Range("K74").Value = 1366.1177
$ws.Range("M74").Value = -492.1177
$ws.Range("H77").Value = 1475.3334
$ws.Range("I77").Value = 1366.1177
$ws.Range("K77").Value = 6830.5885
$ws.Range("M77").Value = -2462.5885
$ws.Range("H98").Value = 30088.5
$ws.Range("J98").Value = 30088.5
$ws.Range("L98").Value = 30088.5
$ws.Range("N98").Value = -36078.5
$ws.Range("H106").Value = 30369.666
$ws.Range("J106").Value = 30369.666
$ws.Range("L106").Value = 30369.666
$ws.Range("N106").Value = -32893.666
$ws.Range("H132").Value = 3459.7942
$ws.Range("I132").Value = 2745.9656
$ws.Range("K132").Value = 8237.8968
$ws.Range("M132").Value = -5707.8968

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1582.75
$ws.Range("I105").Value = 1667.2858
$ws.Range("K105").Value = 1667.2858
$ws.Range("M105").Value = 79.71419999999989
$ws.Range("H134").Value = 1235.3334
$ws.Range("I134").Value = 1235.3334
$ws.Range("K134").Value = 3706.0002
$ws.Range("M134").Value = -1171.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4863.0713
$ws.Range("I31").Value = 2460.818
$ws.Range("J31").Value = 13671.333
$ws.Range("K31").Value = 2460.818
$ws.Range("L31").Value = 13671.333
$ws.Range("M31").Value = -2165.818
$ws.Range("N31").Value = -14261.333
$ws.Range("H34").Value = 4863.0713
$ws.Range("I34").Value = 2460.818
$ws.Range("J34").Value = 13671.333
$ws.Range("K34").Value = 2460.818
$ws.Range("L34").Value = 13671.333
$ws.Range("M34").Value = -2258.818
$ws.Range("N34").Value = -14075.333
$ws.Range("H43").Value = 30328.5
$ws.Range("J43").Value = 30328.5
$ws.Range("L43").Value = 30328.5
$ws.Range("N43").Value = -30696.5
$ws.Range("H57").Value = 78899
$ws.Range("I57").Value = 87500
$ws.Range("J57").Value = 44495
$ws.Range("K57").Value = 87500
$ws.Range("L57").Value = 44495
$ws.Range("M57").Value = -86940
$ws.Range("N57").Value = -45615
$ws.Range("H58").Value = 1317.0625
$ws.Range("I58").Value = 1350.6
$ws.Range("K58").Value = 1350.6
$ws.Range("M58").Value = -1147.6
$ws.Range("H101").Value = 30328.5
$ws.Range("J101").Value = 30328.5
$ws.Range("L101").Value = 30328.5
$ws.Range("N101").Value = -36818.5
$ws.Range("H105").Value = 3386.6667
$ws.Range("I105").Value = 3455
$ws.Range("K105").Value = 3455
$ws.Range("M105").Value = -1708
$ws.Range("H132").Value = 15137.5
$ws.Range("I132").Value = 15137.5
$ws.Range("K132").Value = 45412.5
$ws.Range("M132").Value = -42882.5
$ws.Range("H134").Value = 2630.2173
$ws.Range("I134").Value = 1524.75
$ws.Range("K134").Value = 4574.25
$ws.Range("M134").Value = -2039.25
$ws.Range("H136").Value = 1317.0625
$ws.Range("I136").Value = 1350.6
$ws.Range("K136").Value = 4051.8
$ws.Range("M136").Value = -1501.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4771.5454
$ws.Range("J5").Value = 4931.4443
$ws.Range("L5").Value = 14794.3329
$ws.Range("N5").Value = -15018.3329
$ws.Range("H60").Value = 745.6667
$ws.Range("I60").Value = 511.66666
$ws.Range("J60").Value = 901.6667
$ws.Range("K60").Value = 1534.99998
$ws.Range("L60").Value = 2705.0001
$ws.Range("M60").Value = -1283.99998
$ws.Range("N60").Value = -3207.0001
$ws.Range("H68").Value = 1799.8
$ws.Range("J68").Value = 1666.3334
$ws.Range("L68").Value = 4999.0002
$ws.Range("N68").Value = -6621.0002
$ws.Range("H71").Value = 1799.8
$ws.Range("J71").Value = 1666.3334
$ws.Range("L71").Value = 14997.0006
$ws.Range("N71").Value = -23109.0006
$ws.Range("H123").Value = 1220
$ws.Range("I123").Value = 1220
$ws.Range("K123").Value = 3660
$ws.Range("M123").Value = -1210
$ws.Range("H131").Value = 5950.385
$ws.Range("I131").Value = 12734.777
$ws.Range("J131").Value = 2358.647
$ws.Range("K131").Value = 38204.331
$ws.Range("L131").Value = 7075.941
$ws.Range("M131").Value = -33164.331
$ws.Range("N131").Value = -17155.941
$ws.Range("H135").Value = 4771.5454
$ws.Range("J135").Value = 4931.4443
$ws.Range("L135").Value = 44382.9987
$ws.Range("N135").Value = -49452.9987

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 30671
$ws.Range("J104").Value = 30671
$ws.Range("L104").Value = 30671
$ws.Range("N104").Value = -37659
$ws.Range("H105").Value = 25000
$ws.Range("J105").Value = 25000
$ws.Range("L105").Value = 25000
$ws.Range("N105").Value = -31988
$ws.Range("H134").Value = 45991.168
$ws.Range("J134").Value = 45991.168
$ws.Range("L134").Value = 137973.504
$ws.Range("N134").Value = -143043.504
$ws.Range("H136").Value = 31973.393
$ws.Range("J136").Value = 31973.393
$ws.Range("L136").Value = 95920.179
$ws.Range("N136").Value = -101020.179

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 17900.182
$ws.Range("I93").Value = 2504.9473
$ws.Range("K93").Value = 2504.9473
$ws.Range("M93").Value = -1256.9473
$ws.Range("H103").Value = 27801
$ws.Range("J103").Value = 27801
$ws.Range("L103").Value = 27801
$ws.Range("N103").Value = -30145
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 80305
$ws.Range("J137").Value = 80305
$ws.Range("L137").Value = 80305
$ws.Range("N137").Value = -90505

Write-Host "Applied market price updates across all sheets."